# Updates the "cryptos" price/volume snapshot table (GitHub Actions refresh).
# Note: several Price-column values look like plain numbers (e.g. "597.85"),
# but the sheet stores them as text. Prefixing with a literal leading
# apostrophe ("'" + value) makes Excel keep/enter them as text instead of
# auto-converting to a numeric value, matching the original cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.734.48'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.796.41'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'" + '597.85'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = "'" + '167.71'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').Value = '3.794.03'
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'" + '0.528'
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').Value = "'" + '0.164'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').Value = "'" + '0.461'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').Value = "'" + '0.0000274'
$ws.Range('E13').Value = '  +10.04%  '
$ws.Range('D14').Value = "'" + '36.76'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '4.432.39'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '3.790.57'
$ws.Range('E16').Value = '  -1.95%  '
$ws.Range('D17').Value = '67.844.74'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = "'" + '18.29'
$ws.Range('D19').Value = "'" + '7.43'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = "'" + '10.83'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').Value = "'" + '468.80'
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').Value = "'" + '0.731'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = "'" + '0.0000150'
$ws.Range('E24').Value = '  -8.24%  '
$ws.Range('D25').Value = "'" + '83.28'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').Value = "'" + '2.30'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').Value = "'" + '12.16'
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = "'" + '10.21'
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = "'" + '2.91'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').Value = '3.947.08'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').Value = "'" + '7.70'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = "'" + '2.27'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').Value = "'" + '30.75'
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('D35').Value = "'" + '9.27'
$ws.Range('E35').Value = '  -0.98%  '
$ws.Range('D36').Value = '3.761.83'
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('D37').Value = "'" + '0.106'
$ws.Range('E37').Value = '  +2.07%  '
$ws.Range('D38').Value = "'" + '3.72'
$ws.Range('E38').Value = '  +6.92%  '
$ws.Range('D39').Value = "'" + '5.96'
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = "'" + '0.316'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D45').Value = "'" + '8.75'
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('D46').Value = "'" + '1.96'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').Value = "'" + '46.41'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').Value = "'" + '404.69'
$ws.Range('E48').Value = '  -6.84%  '
$ws.Range('D49').Value = "'" + '0.000285'
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = "'" + '141.86'
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = "'" + '0.0359'
$ws.Range('E51').Value = '  +0.27%  '
